# The workbook's prediction output was regenerated from an updated copy of
# ful-path.csv. The row/column labels and the "prediction" column are
# unchanged — only the numeric score in column B for each sample row is
# refreshed with the newly computed value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 26.711635955412348
$ws.Range("B3").Value = 29.139109776271468
